# Update the "想去人数" (want-to-go count) figures in the "展览" and
# "全部类型" sheets to match the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 4700
    $ws.Range("F3").Value = 2550

    if ($sheetName -eq "展览") {
        $ws.Range("F4").Value = 70
        $ws.Range("F10").Value = 1799
        $ws.Range("F12").Value = 4135
        $ws.Range("F13").Value = 50
        $ws.Range("F14").Value = 285
    }
    elseif ($sheetName -eq "全部类型") {
        $ws.Range("F5").Value = 70
        $ws.Range("F14").Value = 1799
        $ws.Range("F16").Value = 4135
        $ws.Range("F17").Value = 50
        $ws.Range("F18").Value = 285
    }
}
